$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header values (row 1, columns B:E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values (columns B:E)
$ws.Range("B2").Value = 12.02188575120635
$ws.Range("C2").Value = 11.142360010331798
$ws.Range("D2").Value = 12.944287081014441
$ws.Range("E2").Value = 11.992186922912184

# Update row 3 values (columns B:E)
$ws.Range("B3").Value = 10.818102188479116
$ws.Range("C3").Value = 10.097932967693936
$ws.Range("D3").Value = 12.207391533827906
$ws.Range("E3").Value = 11.142287735105581

# Update the selection to match the new active range
$ws.Range("B1:E3").Select()
